$wb = $excel.ActiveWorkbook

# --- NhapSanPham (sheet1): remove the obsolete "ID san pham" column A ---
$ws1 = $wb.Worksheets.Item("NhapSanPham")
$ws1.Columns.Item(1).Delete()

# Manually resize the (now) "Hinh anh" column (col G) - clears bestFit/autosize
$ws1.Columns.Item(7).ColumnWidth = 59.29

# Update the remembered selection on NhapSanPham
$ws1.Range("C9").Select() | Out-Null

# --- NhapKho (sheet3) becomes the active/selected tab ---
$ws3 = $wb.Worksheets.Item("NhapKho")
$ws3.Activate()
